$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3356.8
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3356.8
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3356.8
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3706.8

$ws.Range("H53").Value = 1800.28
$ws.Range("I53").Value = 295.46667
$ws.Range("J53").Value = 4057.5
$ws.Range("K53").Value = 295.46667
$ws.Range("L53").Value = 4057.5
$ws.Range("M53").Value = 341.53333
$ws.Range("N53").Value = -5331.5

$ws.Range("H64").Value = 3049.0356
$ws.Range("I64").Value = 2759.3333
$ws.Range("J64").Value = 3186.2632
$ws.Range("K64").Value = 2759.3333
$ws.Range("L64").Value = 3186.2632
$ws.Range("M64").Value = -2511.3333
$ws.Range("N64").Value = -3682.2632

$ws.Range("H67").Value = 3049.0356
$ws.Range("I67").Value = 2759.3333
$ws.Range("J67").Value = 3186.2632
$ws.Range("K67").Value = 2759.3333
$ws.Range("L67").Value = 3186.2632
$ws.Range("M67").Value = -1901.3333
$ws.Range("N67").Value = -4902.263199999999

$ws.Range("H129").Value = 1356.0975
$ws.Range("I129").Value = 549.4
$ws.Range("J129").Value = 1616.3226
$ws.Range("K129").Value = 1648.2
$ws.Range("L129").Value = 4848.9678
$ws.Range("M129").Value = 3351.8
$ws.Range("N129").Value = -14848.9678

$ws.Range("H138").Value = 1801.4634
$ws.Range("I138").Value = 1771.909
$ws.Range("J138").Value = 1835.6842
$ws.Range("K138").Value = 5315.727000000001
$ws.Range("L138").Value = 5507.0526
$ws.Range("M138").Value = -175.7270000000008
$ws.Range("N138").Value = -15787.0526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7482.8354
$ws.Range("I32").Value = 6331.836
$ws.Range("J32").Value = 20335.666
$ws.Range("K32").Value = 6331.836
$ws.Range("L32").Value = 20335.666
$ws.Range("M32").Value = -6044.836
$ws.Range("N32").Value = -20909.666

$ws.Range("H97").Value = 29872.857
$ws.Range("I97").Value = 34701.668
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 34701.668
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -34205.668
$ws.Range("N97").Value = -1892

$ws.Range("H132").Value = 5363.45
$ws.Range("I132").Value = 4561.3
$ws.Range("J132").Value = 5630.8335
$ws.Range("K132").Value = 13683.9
$ws.Range("L132").Value = 16892.5005
$ws.Range("M132").Value = -11153.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2873.875
$ws.Range("I22").Value = 2998.7144
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 2998.7144
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -2825.7144
$ws.Range("N22").Value = -2346

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H76").Value = 27500
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 27500
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 27500
$ws.Range("N76").Value = -28130

$ws.Range("H79").Value = 27500
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 27500
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 27500
$ws.Range("N79").Value = -29684

$ws.Range("H94").Value = 21430.8
$ws.Range("I94").Value = 26463.625
$ws.Range("J94").Value = 1299.5
$ws.Range("K94").Value = 26463.625
$ws.Range("L94").Value = 1299.5
$ws.Range("M94").Value = -26012.625
$ws.Range("N94").Value = -2201.5

$ws.Range("H105").Value = 1734.4445
$ws.Range("I105").Value = 1355.9
$ws.Range("J105").Value = 2207.625
$ws.Range("K105").Value = 1355.9
$ws.Range("L105").Value = 2207.625
$ws.Range("M105").Value = 391.0999999999999
$ws.Range("N105").Value = -5701.625

$ws.Range("H134").Value = 5008.0933
$ws.Range("I134").Value = 1968.5186
$ws.Range("J134").Value = 10137.375
$ws.Range("K134").Value = 5905.5558
$ws.Range("L134").Value = 30412.125
$ws.Range("M134").Value = -3370.5558
$ws.Range("N134").Value = -35482.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5558039.5
$ws.Range("I31").Value = 2316.6155
$ws.Range("J31").Value = 9806534
$ws.Range("K31").Value = 2316.6155
$ws.Range("L31").Value = 9806534
$ws.Range("M31").Value = -2021.6155
$ws.Range("N31").Value = -9807124

$ws.Range("H34").Value = 5558039.5
$ws.Range("I34").Value = 2316.6155
$ws.Range("J34").Value = 9806534
$ws.Range("K34").Value = 2316.6155
$ws.Range("L34").Value = 9806534
$ws.Range("M34").Value = -2114.6155
$ws.Range("N34").Value = -9806938

$ws.Range("H58").Value = 6821.524
$ws.Range("I58").Value = 3989.1428
$ws.Range("J58").Value = 8237.714
$ws.Range("K58").Value = 3989.1428
$ws.Range("L58").Value = 8237.714
$ws.Range("M58").Value = -3786.1428
$ws.Range("N58").Value = -8643.714

$ws.Range("H62").Value = 11573.182
$ws.Range("I62").Value = 2329.2856
$ws.Range("J62").Value = 27750
$ws.Range("K62").Value = 2329.2856
$ws.Range("L62").Value = 27750
$ws.Range("M62").Value = -1705.2856

$ws.Range("H65").Value = 11573.182
$ws.Range("I65").Value = 2329.2856
$ws.Range("J65").Value = 27750
$ws.Range("K65").Value = 11646.428
$ws.Range("L65").Value = 138750
$ws.Range("M65").Value = -8526.428

$ws.Range("H110").Value = 64700
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 64700
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 64700
$ws.Range("N110").Value = -72880

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H136").Value = 6821.524
$ws.Range("I136").Value = 3989.1428
$ws.Range("J136").Value = 8237.714
$ws.Range("K136").Value = 11967.4284
$ws.Range("L136").Value = 24713.142
$ws.Range("M136").Value = -9417.428400000001
$ws.Range("N136").Value = -29813.142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1861
$ws.Range("I98").Value = 1388.8889
$ws.Range("J98").Value = 2333.111
$ws.Range("K98").Value = 4166.6667
$ws.Range("L98").Value = 6999.333
$ws.Range("M98").Value = -2668.6667
$ws.Range("N98").Value = -9995.332999999999

$ws.Range("H107").Value = 299.45456
$ws.Range("I107").Value = 378.9091
$ws.Range("J107").Value = 220
$ws.Range("K107").Value = 1136.7273
$ws.Range("L107").Value = 660
$ws.Range("M107").Value = 783.2727
$ws.Range("N107").Value = -4500

$ws.Range("H115").Value = 2506.4443
$ws.Range("I115").Value = 2190.875
$ws.Range("J115").Value = 5031
$ws.Range("K115").Value = 6572.625
$ws.Range("L115").Value = 15093
$ws.Range("M115").Value = -5397.625

$ws.Range("H122").Value = 3350.3438
$ws.Range("I122").Value = 1250
$ws.Range("J122").Value = 3490.3667
$ws.Range("K122").Value = 11250
$ws.Range("L122").Value = 31413.3003
$ws.Range("M122").Value = -8800
$ws.Range("N122").Value = -36313.3003

$ws.Range("H132").Value = 1202.75
$ws.Range("I132").Value = 881.4286
$ws.Range("J132").Value = 1452.6666
$ws.Range("K132").Value = 7932.8574
$ws.Range("L132").Value = 13073.9994
$ws.Range("M132").Value = -5402.8574
$ws.Range("N132").Value = -18133.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 710
$ws.Range("I9").Value = 725
$ws.Range("J9").Value = 700
$ws.Range("K9").Value = 725
$ws.Range("L9").Value = 700
$ws.Range("M9").Value = -555
$ws.Range("N9").Value = -1040

$ws.Range("H102").Value = 1694.6786
$ws.Range("I102").Value = 1182.4736
$ws.Range("J102").Value = 2776
$ws.Range("K102").Value = 1182.4736
$ws.Range("L102").Value = 2776
$ws.Range("M102").Value = 439.5264
$ws.Range("N102").Value = -6020

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 759
$ws.Range("I9").Value = 198.75
$ws.Range("J9").Value = 3000
$ws.Range("K9").Value = 198.75
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = 25.25
$ws.Range("N9").Value = -3448

$ws.Range("H30").Value = 70018
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 70018
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 70018
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -70234

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 10000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 10000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 10000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -11872

$ws.Range("H78").Value = 10000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 10000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 30000
$ws.Range("M78").ClearContents()
